$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $val) {
    $c = $ws.Range($cell)
    $origStyle = $c.Style
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = $origStyle
}

$ws.Range('D2').Value = '26.701.63'
$ws.Range('E2').Value = '  -1.89%  '
$ws.Range('D3').Value = '1.681.26'
$ws.Range('E3').Value = '  -1.34%  '
$ws.Range('E4').Value = '  +0.56%  '
Set-TextValue 'D5' '220.98'
$ws.Range('E5').Value = '  -0.79%  '
Set-TextValue 'D6' '0.5227'
$ws.Range('E6').Value = '  -1.28%  '
Set-TextValue 'D7' '1.008'
$ws.Range('E7').Value = '  +0.50%  '
Set-TextValue 'D8' '0.06548'
$ws.Range('E8').Value = '  +0.02%  '
Set-TextValue 'D9' '0.2589'
$ws.Range('E9').Value = '  -1.96%  '
Set-TextValue 'D10' '20.24'
$ws.Range('E10').Value = '  -2.62%  '
Set-TextValue 'D11' '0.07713'
$ws.Range('E11').Value = '  +1.06%  '
$ws.Range('D12').Value = '1.723.76'
$ws.Range('E12').Value = '  +1.10%  '
$ws.Range('D13').Value = '1.920.35'
$ws.Range('E13').Value = '  -1.01%  '
$ws.Range('E14').Value = '  -4.47%  '
Set-TextValue 'D15' '0.5633'
$ws.Range('E15').Value = '  -1.56%  '
$ws.Range('D16').Value = '0.0₅8072'
$ws.Range('E16').Value = '  -1.13%  '
Set-TextValue 'D17' '65.75'
$ws.Range('E17').Value = '  -2.46%  '
$ws.Range('D18').Value = '26.773.77'
$ws.Range('E18').Value = '  -1.61%  '
Set-TextValue 'D19' '215.02'
$ws.Range('E19').Value = '  -0.34%  '
$ws.Range('E20').Value = '  +0.32%  '
Set-TextValue 'D21' '4.525'
$ws.Range('E21').Value = '  -2.88%  '
Set-TextValue 'D22' '10.19'
$ws.Range('E22').Value = '  -2.28%  '
Set-TextValue 'D23' '5.933'
$ws.Range('E23').Value = '  -0.40%  '
Set-TextValue 'D24' '1.008'
$ws.Range('E24').Value = '  +0.51%  '
Set-TextValue 'D25' '143.88'
$ws.Range('E25').Value = '  +1.17%  '
Set-TextValue 'D26' '1.731'
$ws.Range('E26').Value = '  -0.89%  '
Set-TextValue 'D27' '0.1182'
$ws.Range('E27').Value = '  -2.84%  '
Set-TextValue 'D28' '7.079'
$ws.Range('E28').Value = '  -2.29%  '
Set-TextValue 'D29' '15.93'
$ws.Range('E29').Value = '  -1.97%  '
Set-TextValue 'D30' '0.05274'
$ws.Range('E30').Value = '  -1.54%  '
Set-TextValue 'D31' '1.277'
Set-TextValue 'D32' '3.383'
$ws.Range('E32').Value = '  -3.20%  '
Set-TextValue 'D33' '3.265'
$ws.Range('E33').Value = '  -4.19%  '
$ws.Range('E34').Value = '  -1.77%  '
Set-TextValue 'D35' '2.783'
$ws.Range('E35').Value = '  -3.18%  '
Set-TextValue 'D36' '2.394'
$ws.Range('E36').Value = '  -1.03%  '
Set-TextValue 'D37' '0.9333'
$ws.Range('E37').Value = '  -1.27%  '
Set-TextValue 'D38' '0.5831'
$ws.Range('E38').Value = '  -0.08%  '
$ws.Range('D39').Value = '1.169.00'
$ws.Range('E39').Value = '  +12.59%  '
Set-TextValue 'D40' '0.01618'
$ws.Range('E40').Value = '  -0.60%  '
$ws.Range('E41').Value = '  +0.54%  '
Set-TextValue 'D42' '5.727'
$ws.Range('E42').Value = '  -2.38%  '
Set-TextValue 'D43' '0.8326'
$ws.Range('E43').Value = '  -0.59%  '
Set-TextValue 'D44' '100.01'
$ws.Range('E44').Value = '  -0.86%  '
$ws.Range('D45').Value = '1.828.50'
$ws.Range('E45').Value = '  -1.02%  '
$ws.Range('D46').Value = '0.0₈110'
$ws.Range('E46').Value = '  -3.18%  '
$ws.Range('B47').Value = 'Aave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue 'D47' '56.26'
$ws.Range('E47').Value = '  -2.78%  '
$ws.Range('B48').Value = 'Mantle'
$ws.Range('C48').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue 'D48' '0.4502'
$ws.Range('E48').Value = '  +0.35%  '
Set-TextValue 'D49' '1.005'
$ws.Range('E49').Value = '  +0.27%  '
Set-TextValue 'D50' '8.043'
$ws.Range('E50').Value = '  -0.30%  '
Set-TextValue 'D51' '0.05185'
$ws.Range('E51').Value = '  -1.03%  '

Write-Output "Applied 98 cell updates"
